# FlixelRL-536: split the "treasure chest" message into two rows and
# add a new row for the image-display (IMAGE) command text that used to
# share a cell with "<br>" markup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Shorten the existing message in B25 to just the first sentence.
$ws.Range("B25").Value2 = "女の子はドキドキして宝箱を開きます。"

# 2) Adjust row heights that no longer need the extra space reserved for
#    the "<br>" wrapped text.
$ws.Rows.Item(25).RowHeight = 20
$ws.Rows.Item(28).RowHeight = 20
$ws.Rows.Item(31).RowHeight = 20

# 3) Insert a brand-new row 33 carrying the second sentence that used to
#    live inside B25, copying the formatting from the row above it (row
#    32) so the new cells match the surrounding style/border/font.
$ws.Range("A32:B32").Copy() | Out-Null
$ws.Range("A33:B33").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A33").Value2 = 31
$ws.Range("B33").Value2 = "しかし、中にあるのは小さなマタタビの木でした。"
$ws.Rows.Item(33).RowHeight = 20
